$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 37 ("Sandia", fecha 2022-12-05,
# calidad Tercera), pushing the existing rows 37-49 down to 38-50.
$ws.Rows.Item(37).Insert()

$ws.Range("A37").Value = 1
$ws.Range("B37").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C37").Value = "Arica y Parinacota"
$ws.Range("D37").Value = 44900
$ws.Range("E37").Value = 15
$ws.Range("F37").Value = 100112028
$ws.Range("G37").Value = "Sandia"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Tercera"
$ws.Range("J37").Value = 600
$ws.Range("K37").Value = 480
$ws.Range("L37").Value = 500
$ws.Range("M37").Value = 490
$ws.Range("N37").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O37").Value = "Perú"
$ws.Range("P37").Value = 490
$ws.Range("Q37").Value = 1
$ws.Range("R37").Value = "Hortaliza"
